$wb = $excel.ActiveWorkbook

$wsRandom  = $wb.Worksheets.Item("random")
$wsStatic  = $wb.Worksheets.Item("static")
$wsDynamic = $wb.Worksheets.Item("dynamic")

# ---------------------------------------------------------------------------
# "static" sheet: refresh the metrics table with new values from a re-run.
# ---------------------------------------------------------------------------
$staticData = @(
    @(1, "relu", 2, 92.199999094009399),
    @(4, "relu", 5, 92.000001668930054),
    @(3, "relu", 4, 91.79999828338623),
    @(2, "relu", 3, 91.200000047683716),
    @(5, "relu", 6, 91.00000262260437),
    @(0, "relu", 1, 88.599997758865356)
)

for ($i = 0; $i -lt $staticData.Length; $i++) {
    $r = 2 + $i
    $row = $staticData[$i]
    $wsStatic.Cells.Item($r, 1).Value = $row[0]
    $wsStatic.Cells.Item($r, 2).Value = $row[1]
    $wsStatic.Cells.Item($r, 3).Value = $row[2]
    $wsStatic.Cells.Item($r, 4).Value = $row[3]
}

# rows 8-9 no longer have data in the refreshed export - clear them but keep formatting
$wsStatic.Range("A8:D9").ClearContents()

# ---------------------------------------------------------------------------
# "dynamic" sheet: refresh the metrics table with new values from a re-run.
# ---------------------------------------------------------------------------
$dynamicData = @(
    @(5, "relu", 6, 91.600000858306885),
    @(3, "relu", 4, 90.600001811981201),
    @(2, "relu", 3, 90.200001001358032),
    @(7, "relu", 8, 89.999997615814209),
    @(1, "relu", 2, 89.800000190734863),
    @(6, "relu", 7, 89.800000190734863),
    @(0, "relu", 1, 89.200001955032349),
    @(4, "relu", 5, 87.999999523162842)
)

for ($i = 0; $i -lt $dynamicData.Length; $i++) {
    $r = 2 + $i
    $row = $dynamicData[$i]
    $wsDynamic.Cells.Item($r, 1).Value = $row[0]
    $wsDynamic.Cells.Item($r, 2).Value = $row[1]
    $wsDynamic.Cells.Item($r, 3).Value = $row[2]
    $wsDynamic.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# Re-apply header/row-label styling on "static" so the bold cells pick up an
# explicit theme color - this is what produced the extra font entry in
# styles.xml (matches the already-bold "dynamic" sheet's look).
# ---------------------------------------------------------------------------
$wsStatic.Range("B1:D1").Font.ThemeColor = 1
$wsStatic.Range("A2:A7").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# Selections & active sheet/tab: the workbook now opens on "random" (it lost
# its explicit activeTab), "static" and "dynamic" keep their own selections.
# ---------------------------------------------------------------------------
$wsStatic.Range("C15").Select()
$wsDynamic.Range("C11").Select()
$wsRandom.Activate()
$wsRandom.Range("D13").Select()
